# Progress Tracking.xlsx - "update status for 21 Nov"
# Fills in the Daily Activity log for 18-22 Nov 2020 (rows 26-32), extends the
# running total to row 33, merges the new combined Nov-21 entry (rows 29-30),
# and flips the "13 Networking" completion flag to Yes on the Topics sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# Grab copies of formatting we will need *after* we start overwriting rows,
# by reading it from rows that are not touched by this edit.
# ---------------------------------------------------------------------------

# Row 26 ("No Work" day 18-Nov). Date cell (B26) already has the right date &
# format, only the rest of the row needs to be filled in using the same
# look as the existing "No Work" row (row 18).
$ws.Range("C18:F18").Copy()
$ws.Range("C26:F26").PasteSpecial($xlPasteFormats)
$ws.Range("C26").Value = "No Work"
$ws.Range("D26").Value = "No Work"
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 0

# Row 27 ("No Work" day 19-Nov) - same look as row 18, with a new date.
$ws.Range("B18").Copy()
$ws.Range("B27").PasteSpecial($xlPasteFormats)
$ws.Range("B27").Value = 44154
$ws.Range("C18:F18").Copy()
$ws.Range("C27:F27").PasteSpecial($xlPasteFormats)
$ws.Range("C27").Value = "No Work"
$ws.Range("D27").Value = "No Work"
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 0

# Row 28 (13 Networking, 3 videos 90-92, 20-Nov) - normal single-day row,
# same look as row 21.
$ws.Range("B18").Copy()
$ws.Range("B28").PasteSpecial($xlPasteFormats)
$ws.Range("B28").Value = 44155
$ws.Range("C21:F21").Copy()
$ws.Range("C28:F28").PasteSpecial($xlPasteFormats)
$ws.Range("C28").Value = "13 Networking"
$ws.Range("D28").Value = "3 videos, 90-92"
$ws.Range("E28").Value = 43
$ws.Range("F28").Value = 43

# Row 29-30 (21-Nov, merged): 13 Networking / create project on Networking
# demo, then 15 Best Practices BONUS / download ebook. Uses the same
# "two day merged" look as rows 16-17.
$ws.Range("B16").Copy()
$ws.Range("B29").PasteSpecial($xlPasteFormats)
$ws.Range("B17").Copy()
$ws.Range("B30").PasteSpecial($xlPasteFormats)
$ws.Range("F16").Copy()
$ws.Range("F29").PasteSpecial($xlPasteFormats)
$ws.Range("F17").Copy()
$ws.Range("F30").PasteSpecial($xlPasteFormats)
$ws.Range("C21:E21").Copy()
$ws.Range("C29:E29").PasteSpecial($xlPasteFormats)
$ws.Range("C21:E21").Copy()
$ws.Range("C30:E30").PasteSpecial($xlPasteFormats)

$ws.Range("B29").Value = 44156
$ws.Range("C29").Value = "13 Networking"
$ws.Range("D29").Value = "create project on Networking demo"
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 0

$ws.Range("C30").Value = "15 Best Practices  BONUS"
$ws.Range("D30").Value = "download ebook & pdf from qmlbook.github.io"
$ws.Range("E30").Value = 0

$ws.Range("B29:B30").Merge()
$ws.Range("F29:F30").Merge()

# Row 31-32: new blank placeholder rows (same look the old placeholder rows
# 27-28 used to have), row 31 carries the next date (22-Nov).
$ws.Range("B18").Copy()
$ws.Range("B31").PasteSpecial($xlPasteFormats)
$ws.Range("B31").Value = 44157

$ws.Range("C16:E16").Copy()
$ws.Range("C31:E31").PasteSpecial($xlPasteFormats)
$ws.Range("F21").Copy()
$ws.Range("F31").PasteSpecial($xlPasteFormats)

$ws.Range("C32").Value = ""
$ws.Range("C16:E16").Copy()
$ws.Range("C32:E32").PasteSpecial($xlPasteFormats)
$ws.Range("F8").Copy()
$ws.Range("F32").PasteSpecial($xlPasteFormats)

$ws.Range("C31").Value = ""
$ws.Range("D31").Value = ""
$ws.Range("E31").Value = ""
$ws.Range("C32").Value = ""
$ws.Range("D32").Value = ""
$ws.Range("E32").Value = ""

# Row 33: running total now covers F3:F32, using the same look the old
# F29 total cell had.
$ws.Range("F29").Copy()
$ws.Range("F33").PasteSpecial($xlPasteFormats)
$ws.Range("F33").Formula = "=SUM(F3:F32)"

$excel.CutCopyMode = $false

# Update the view to roughly match where the edit was made.
$ws.Range("A18").Select()
$ws.Range("H31").Select()

# ---------------------------------------------------------------------------
# Topics sheet: "13 Networking" is now complete.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("E15").Value = "Yes"
